# Auto-generated edit script: updates cryptocurrency price/volume data
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.567.02"
$ws.Range("E2").Value = "  +2.39%  "
# Row 3
$ws.Range("D3").Value = "2.289.96"
$ws.Range("E3").Value = "  +1.58%  "
# Row 4
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "157.56"
$ws.Range("E5").Value = "  +15,629.64%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.75"
$ws.Range("E6").Value = "  +1.26%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "97.00"
$ws.Range("E7").Value = "  +6.44%  "
# Row 8
$ws.Range("E8").Value = "  +0.70%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.497"
$ws.Range("E10").Value = "  +4.04%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "36.38"
$ws.Range("E11").Value = "  +13.59%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0806"
$ws.Range("E12").Value = "  +1.41%  "
# Row 13
$ws.Range("E13").Value = "  -1.70%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.75"
$ws.Range("E14").Value = "  +2.84%  "
# Row 15
$ws.Range("D15").Value = "2.643.47"
$ws.Range("E15").Value = "  +1.54%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.61"
$ws.Range("E16").Value = "  +3.21%  "
# Row 17
$ws.Range("D17").Value = "2.311.13"
$ws.Range("E17").Value = "  +2.30%  "
# Row 18
$ws.Range("E18").Value = "  +6.32%  "
# Row 19
$ws.Range("D19").Value = "42.430.83"
$ws.Range("E19").Value = "  +2.22%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.80"
$ws.Range("E20").Value = "  +4.17%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0921"
$ws.Range("E21").Value = "  +2.21%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.02"
$ws.Range("E22").Value = "  +2.48%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.95"
$ws.Range("E23").Value = "  +2.16%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.48"
$ws.Range("E24").Value = "  +1.55%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.60"
$ws.Range("E25").Value = "  +1.50%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.96"
$ws.Range("E26").Value = "  +2.91%  "
# Row 27
$ws.Range("E27").Value = "  +0.08%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.07"
$ws.Range("E28").Value = "  +0.92%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.33"
$ws.Range("E29").Value = "  +6.57%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.62"
$ws.Range("E30").Value = "  +1.61%  "
# Row 31
$ws.Range("E31").Value = "  -1.00%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.03"
$ws.Range("E32").Value = "  +0.62%  "
# Row 33
$ws.Range("E33").Value = "  +4.00%  "
# Row 34
$ws.Range("E34").Value = "  -0.05%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0754"
$ws.Range("E35").Value = "  +1.98%  "
# Row 36
$ws.Range("E36").Value = "  +3.20%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.37"
$ws.Range("E37").Value = "  +4.91%  "
# Row 38
$ws.Range("E38").Value = "  +4.98%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("E39").Value = "  +5.38%  "
# Row 41
$ws.Range("E41").Value = "  -0.07%  "
# Row 42
$ws.Range("E42").Value = "  +7.97%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.33"
$ws.Range("E43").Value = "  +14.63%  "
# Row 44
$ws.Range("D44").Value = "2.005.40"
$ws.Range("E44").Value = "  -2.18%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.28"
$ws.Range("E45").Value = "  -1.28%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0284"
$ws.Range("E46").Value = "  +2.53%  "
# Row 47
$ws.Range("E47").Value = "  +6.51%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.25"
$ws.Range("E48").Value = "  +0.79%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.24"
$ws.Range("E49").Value = "  +5.54%  "
# Row 50
$ws.Range("E50").Value = "  +2.20%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.89"
$ws.Range("E51").Value = "  +0.75%  "
